$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: theta_se values (in parentheses)
$ws.Range("B4").Value = "(1.97)"
$ws.Range("C4").Value = "(0.89)"
$ws.Range("D4").Value = "(0.45)"
$ws.Range("E4").Value = "(1.36)"
$ws.Range("F4").Value = "(0.54)"
$ws.Range("G4").Value = "(0.35)"
$ws.Range("H4").Value = "(0.3)"
$ws.Range("I4").Value = "(0.12)"
$ws.Range("J4").Value = "(0.96)"

# Row 6: lambda_se values (in parentheses)
$ws.Range("B6").Value = "(1.52)"
$ws.Range("C6").Value = "(0.86)"
$ws.Range("D6").Value = "(0.23)"
$ws.Range("E6").Value = "(0.88)"
$ws.Range("F6").Value = "(0.68)"
$ws.Range("G6").Value = "(0.28)"
$ws.Range("H6").Value = "(0.0)"
$ws.Range("I6").Value = "(0.91)"
$ws.Range("J6").Value = "(0.05)"
